$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("S","T","U","V","W","X","Y")

for ($r = 2; $r -le 81; $r++) {
    $ws.Range("C$r").Value = 45208

    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        $f = $cell.Formula
        if ($f -and ($f -like "*Logging_HEBY*")) {
            $cell.Formula = ($f -replace "Logging_HEBY", "Logging_0331")
        }
    }
}
